# update as of 5/10
# Correct the Date column (F) for rows 880-952: these entries were
# recorded as 9/1/2022 (serial 44848) but should be 8/31/2022 (serial 44847).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F880:F952").Value = 44847

# Reset the view back to the top-left / A1, clearing the scrolled-down
# "topLeftCell=A929" / "selection=H936" state that was left over from
# editing near the bottom of the sheet.
$ws.Range("A1").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
